# "in out paras absorber" - clear the placeholder in/out parameter values
# that were filled across the stream table (columns C:M, rows 4-23) on the
# STREAMS sheet, then leave that sheet as the active/selected one with the
# cleared range selected (mirrors the manual edit captured in the diff).

$wb = $excel.ActiveWorkbook

$streams = $wb.Worksheets.Item("STREAMS")

# Clear all the per-stream numeric values that were entered as placeholders.
# ClearContents removes the cached value but keeps any explicit cell
# formatting (e.g. K6/L6 keep their scientific-notation style as an empty,
# styled cell instead of disappearing entirely).
$streams.Range("C4:M23").ClearContents() | Out-Null

# The STREAMS sheet becomes the active tab/selection (previously it was
# "UNIT OPERATIONS"), with the cleared block left selected.
$streams.Activate() | Out-Null
$streams.Range("C4:M23").Select() | Out-Null
